$d = $word.ActiveDocument

# 1. Title heading and the matching bold reprise near the end (both occurrences
#    share identical text, so a single whole-document ReplaceAll handles both).
$d.Content.Find.Execute(
    "Play Million Coins Respin for Free - Classic Frut Machine Style",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Million Coins Respins Free - Classic Slot Gameplay", 2)

# 2. "What we like" bullet list.
#    Processed in reverse (last bullet first) because the new text for the
#    first bullet ("Traditional fruit machine style") is a substring of the
#    old text of the last bullet ("Traditional fruit machine style appeals to
#    classic slot fans"). Replacing the last bullet first avoids the new text
#    being erroneously matched/altered inside the not-yet-updated last bullet.
$d.Content.Find.Execute(
    "Traditional fruit machine style appeals to classic slot fans",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Classic slot gaming experience", 2)

$d.Content.Find.Execute(
    "Simple and easy-to-understand gameplay",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Manual Respin feature for big wins", 2)

$d.Content.Find.Execute(
    "Manual Respin feature can lead to big wins",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Six jackpot prizes", 2)

$d.Content.Find.Execute(
    "Offers six jackpot prizes with high payouts",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Traditional fruit machine style", 2)

# 3. "What we don't like" bullet.
$d.Content.Find.Execute(
    "Chances of winning big with Respin feature are lower",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lower chances of winning with Respin feature", 2)

# 4. Meta description italic paragraph near the end.
$d.Content.Find.Execute(
    "Read our review of Million Coins Respin, a classic fruit machine style slot game from iSoftBet. Play for free and win big with six jackpot prizes and a Respin feature.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Million Coins Respins and play for free. Experience classic slot gaming with six jackpot prizes and a Respin feature.", 2)
